$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.51%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.25%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.711"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.31%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08375"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.83%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.024"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.797"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.93%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.546"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.40%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9271"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.00%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1296"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.54%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1971"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.26%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09544"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.77%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03901"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6.66%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1062"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.76%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001307"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.06%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006150"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.75%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.442"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.44%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.11%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.238"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.35%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.32%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04415"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001255"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.40%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004370"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.12%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02804"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.73%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.77%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007952"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.47%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.33%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009222"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002142"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.75%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01107"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.08%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007013"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.65%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.40%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003513"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "15.48%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.33%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.40%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.40%"
